# Add a new "2022-Q3" sheet right after "总计" / before "2022-Q2",
# fill it with fund-holding detail data, and update the "总计"
# (summary) sheet's table with the new quarter's row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Insert the new worksheet in the right tab position
# ---------------------------------------------------------------
$sheetQ2ref = $wb.Worksheets.Item("2022-Q2")
$newSheet = $wb.Worksheets.Add($sheetQ2ref)
$newSheet.Name = "2022-Q3"

# Match the page-setup margins used by its sibling quarter sheets
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# Re-fetch sheet objects by name (fresh references) so cross-sheet
# Range.Copy formatting actually takes effect.
$template = $wb.Worksheets.Item("2022-Q2")
$ws = $wb.Worksheets.Item("2022-Q3")

# Seed the exact same look (bold/centered/bordered header + index
# column, plain data cells) used by the sibling quarter sheets.
$template.Range("B1:H1").Copy($ws.Range("B1:H1"))
$template.Range("A2:A6").Copy($ws.Range("A2:A6"))
$template.Range("B2:G6").Copy($ws.Range("B2:G6"))
$template.Range("H2:H6").Copy($ws.Range("H2:H6"))

# ---------------------------------------------------------------
# 2. Header row text
# ---------------------------------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, 2 + $i).Value = $headers[$i]
}

# ---------------------------------------------------------------
# 3. Data rows
# ---------------------------------------------------------------
$rows = @(
    @("005777", "广发科技动力股票",        "16.94", "84.51", "6.64", "1.1248", 4),
    @("007731", "民生加银持续成长混合A",    "3.22",  "94.57", "9.36", "0.3014", 2),
    @("007732", "民生加银持续成长混合C",    "1.89",  "94.57", "9.36", "0.1769", 2),
    @("005310", "广发电子信息传媒股票A",    "1.55",  "89.36", "2.89", "0.0448", 8),
    @("010236", "广发电子信息传媒股票C",    "0.13",  "89.36", "2.89", "0.0038", 8)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $r + 2
    $data = $rows[$r]

    # Index column (already numbered 0..4 by the formatting copy above)
    $ws.Cells.Item($row, 1).Value = $r

    # B..G are stored as text (keeps leading zeros / exact decimal text)
    for ($c = 0; $c -lt 6; $c++) {
        $cell = $ws.Cells.Item($row, 2 + $c)
        $cell.NumberFormat = "@"
        $cell.Value = $data[$c]
    }

    # H is numeric
    $ws.Cells.Item($row, 8).Value = $data[6]
}

# ---------------------------------------------------------------
# 4. Update the "总计" (summary) sheet: insert the 2022-Q3 row at
#    the top of the table and shift the rest down.
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Range("A6:D6").Copy($summary.Range("A7:D7"))

$summary.Range("B7").Value = "2021-Q2"
$summary.Range("C7").Value = 16
$summary.Range("D7").Value = 1.54

$summary.Range("B6").Value = "2021-Q3"
$summary.Range("C6").Value = 12
$summary.Range("D6").Value = 3.48

$summary.Range("B5").Value = "2021-Q4"
$summary.Range("C5").Value = 9
$summary.Range("D5").Value = 5.54

$summary.Range("B4").Value = "2022-Q1"
$summary.Range("C4").Value = 10
$summary.Range("D4").Value = 2.11

$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 9
$summary.Range("D3").Value = 2.1

$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 5
$summary.Range("D2").Value = 1.65

$summary.Range("A2").Value = 0
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5

# ---------------------------------------------------------------
# 5. Restore the originally active sheet
# ---------------------------------------------------------------
$wb.Worksheets.Item("总计").Activate()
